$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 120 was a blank "spacer" row (only formatted, no values). A new time
# entry (2014-04-09, 0:68055... -> 0.71527...) is inserted here, extending
# the A:G data table by one row.
# ---------------------------------------------------------------------------
$ws.Range("A120").Value = 2014
$ws.Range("B120").Value = 4
$ws.Range("C120").Value = 9
$ws.Range("D120").Value = 0.68055555555555547
$ws.Range("E120").Value = 0.71527777777777779
$ws.Range("F120").Formula = "=(E120-D120)*24*60"
$ws.Range("G120").Formula = "=F120/60"

# ---------------------------------------------------------------------------
# Everything below the data table shifts down by one row: the blank spacer
# row moves from 120 -> 121, and the three summary rows move from
# 121/122/123 -> 122/123/124. A new summary row is appended at the end.
# ---------------------------------------------------------------------------

# Row 121: new blank spacer row (same look as the old row 120).
$ws.Range("E121").ClearContents()
$ws.Range("F121").ClearContents()
$ws.Range("E121").ClearFormats()
$ws.Range("D121").NumberFormat = "hh:mm;@"
$ws.Range("E121").NumberFormat = "hh:mm;@"
$ws.Range("F121").NumberFormat = "0"

# Row 122: "sum [min]" (previously row 121), sums now go through row 121.
$ws.Range("E122").Value = "sum [min]"
$ws.Range("F122").Formula = "=SUM(F2:F121)"
$ws.Range("F122").NumberFormat = "0"

# Row 123: "sum [h]" (previously row 122).
$ws.Range("E123").Value = "sum [h]"
$ws.Range("F123").Formula = "=F122/60"
$ws.Range("F123").NumberFormat = "0.00"

# Row 124: "sum [working weeks]" (previously row 123, now appended new).
$ws.Range("E124").Value = "sum [working weeks]"
$ws.Range("E124").HorizontalAlignment = -4152
$ws.Range("F124").Formula = "=F123/38.5"
$ws.Range("F124").NumberFormat = "0.00"

# ---------------------------------------------------------------------------
# View bookkeeping: keep the same selected cell but scroll so the new last
# rows of the table are visible (topLeftCell A106 -> A103).
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 103
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I120").Select()
